$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new region/city data rows 180-219 (column B = city/region name, column H = empire name on header rows)
$ws.Cells.Item(180, 2).Value = 'Хэйан-кё'
$ws.Cells.Item(180, 8).Value = 'Япония'
$ws.Cells.Item(181, 2).Value = 'Нагаока'
$ws.Cells.Item(182, 2).Value = 'Нара'
$ws.Cells.Item(183, 2).Value = 'Осака'
$ws.Cells.Item(184, 2).Value = 'Киото'
$ws.Cells.Item(185, 2).Value = 'Хиросима'
$ws.Cells.Item(186, 2).Value = 'Кобе'
$ws.Cells.Item(187, 2).Value = 'Токио'
$ws.Cells.Item(188, 2).Value = 'Дамаск'
$ws.Cells.Item(188, 8).Value = 'Аббасидский халифат'
$ws.Cells.Item(189, 2).Value = 'Багдад'
$ws.Cells.Item(190, 2).Value = 'Самарра'
$ws.Cells.Item(191, 2).Value = 'Александрия'
$ws.Cells.Item(192, 2).Value = 'Медина'
$ws.Cells.Item(193, 2).Value = 'Алеппо'
$ws.Cells.Item(194, 2).Value = 'Басра'
$ws.Cells.Item(195, 2).Value = 'Мекка'
$ws.Cells.Item(196, 2).Value = 'Саксония'
$ws.Cells.Item(196, 8).Value = 'Восточно-Франкское'
$ws.Cells.Item(197, 2).Value = 'Тюрингия'
$ws.Cells.Item(198, 2).Value = 'Бавария'
$ws.Cells.Item(199, 2).Value = 'Алемания'
$ws.Cells.Item(200, 2).Value = 'Франкфурт-на-Майне'
$ws.Cells.Item(201, 2).Value = 'Бремен'
$ws.Cells.Item(202, 2).Value = 'Падеборн'
$ws.Cells.Item(203, 2).Value = 'Майнц'
$ws.Cells.Item(204, 2).Value = 'Регенсбург'
$ws.Cells.Item(205, 2).Value = 'Кур'
$ws.Cells.Item(206, 2).Value = 'Гамбург'
$ws.Cells.Item(207, 2).Value = 'Нанси'
$ws.Cells.Item(207, 8).Value = 'Лотарингия'
$ws.Cells.Item(208, 2).Value = 'Кельн'
$ws.Cells.Item(209, 2).Value = 'Льеж'
$ws.Cells.Item(210, 2).Value = 'Верден'
$ws.Cells.Item(211, 2).Value = 'Безансон'
$ws.Cells.Item(212, 2).Value = 'Ахен'
$ws.Cells.Item(213, 2).Value = 'Нант'
$ws.Cells.Item(213, 8).Value = 'Западно-Франкское'
$ws.Cells.Item(214, 2).Value = 'Тур'
$ws.Cells.Item(215, 2).Value = 'Париж'
$ws.Cells.Item(216, 2).Value = 'Нарбонн'
$ws.Cells.Item(217, 2).Value = 'Бордо'
$ws.Cells.Item(218, 2).Value = 'Барселона'
$ws.Cells.Item(219, 2).Value = 'Реймс'

# Update selection to match final cursor position (B220) as in the authored edit
$ws.Range("B220").Select()
